# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns.
# Values are stored as text in the sheet; a leading "'" forces Excel to
# keep digit-and-dot strings (e.g. "239.70") as text instead of coercing
# them to numbers, matching the original inline-string cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.747.25"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "1.719.44"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "'239.70"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").Value = "'0.4747"
$ws.Range("E7").Value = "  -2.43%  "
$ws.Range("D8").Value = "'0.2541"
$ws.Range("E8").Value = "  -1.70%  "
$ws.Range("D9").Value = "'0.06098"
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("D10").Value = "1.719.32"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("D11").Value = "'15.82"
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("D12").Value = "'0.06872"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").Value = "'0.5940"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").Value = "'4.396"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("D15").Value = "'76.22"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "26.604.94"
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").Value = "'0.000007034"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").Value = "'11.22"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").Value = "1.944.48"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'4.359"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").Value = "'8.298"
$ws.Range("E23").Value = "  -1.79%  "
$ws.Range("D24").Value = "'5.020"
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("D25").Value = "'140.12"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("D26").Value = "'15.10"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("D27").Value = "'1.770"
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("D28").Value = "'105.88"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").Value = "'1.367"
$ws.Range("E29").Value = "  -2.94%  "
$ws.Range("D30").Value = "'3.910"
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("D31").Value = "'0.07847"
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("D32").Value = "'3.607"
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").Value = "'0.04514"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("D34").Value = "'2.616"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "'0.9884"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").Value = "'0.6075"
$ws.Range("E36").Value = "  -1.77%  "
$ws.Range("D37").Value = "'0.9127"
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("D38").Value = "'2.478"
$ws.Range("E38").Value = "  +4.33%  "
$ws.Range("D39").Value = "'1.953"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").Value = "'1.001"
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("D41").Value = "'5.705"
$ws.Range("E41").Value = "  +4.87%  "
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").Value = "'100.08"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").Value = "'0.3774"
$ws.Range("D45").Value = "'6.683"
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").Value = "'0.05347"
$ws.Range("D48").Value = "'7.810"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").Value = "'29.51"
$ws.Range("E49").Value = "  -2.83%  "
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("E51").Value = "  +0.34%  "
